$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'98.831.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.11%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.350.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.61%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'260.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.40%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'650.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.31%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +10.68%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.468"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +19.10%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +22.64%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'3.347.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.63%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +5.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'43.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +21.28%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000270"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +8.66%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'99.523.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.76%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.987.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +1.20%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.352.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.41%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +21.32%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'16.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +10.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'539.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +8.93%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.51%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +8.98%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.0000213"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.60%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.434"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +53.93%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'102.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +14.48%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'6.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +9.45%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'12.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.27%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.527.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.02%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +10.59%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.15%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'11.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +14.84%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.193"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.64%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.00%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'29.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.74%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.536"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +15.76%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'7.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +5.46%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.34%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'519.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.27%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.63%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +4.04%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -2.36%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0426"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +29.98%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'dogwifhat"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'3.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.78%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'ARBITRUM"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.827"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.33%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.61%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +19.24%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'5.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +9.05%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'164.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.25%  "
$ws.Range("E51").Style = "Normal"
